$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = "Picado"
$ws.Range("D2").Value = "ar{kgnaer{hk"

$ws.Range("C3").Value = "Normal"
$ws.Range("D3").Value = "Sin modificaciones"
